# 02-TC-Customers.xlsx update: add pass/fail/blocked/total summary table on
# Sheet1 (I2:M4) and the supporting COUNTIF/SUM breakdown table on Sheet2
# (A6:E17), per the "Updated all test documents added test percentage"
# commit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1 — summary header (row 2) + percentage row (rows 3-4, merged)
# ---------------------------------------------------------------------

# Row 2 labels
$ws1.Range("I2").Value = "Pass"
$ws1.Range("J2").Value = "Fail"
$ws1.Range("K2").Value = "NOT RUN"
$ws1.Range("L2").Value = "Blocked"
$ws1.Range("M2").Value = "Total"

# Row 2 fills (green / red / none / purple / none)
$ws1.Range("I2").Interior.Color = 5296274
$ws1.Range("J2").Interior.Color = 255
$ws1.Range("L2").Interior.Color = 10498160

# Row 2 alignment + thin box border
$hdrRange = $ws1.Range("I2:M2")
$hdrRange.HorizontalAlignment = -4108
$hdrRange.Borders.LineStyle = 1
$hdrRange.Borders.Weight = 2

# Row 3 formulas (percentage of each bucket out of the Sheet2 grand total)
$ws1.Range("I3").Formula = "=AVERAGE(Sheet2!B17/85)"
$ws1.Range("J3").Formula = "=AVERAGE(Sheet2!C17/85)"
$ws1.Range("K3").Formula = "=AVERAGE(Sheet2!D17/85)"
$ws1.Range("L3").Formula = "=AVERAGE(Sheet2!E17/85)"
$ws1.Range("M3").Formula = "=SUM(I3:L4)"

# Row 4 — blank cells under row 3 (merged I3:I4 etc.), same formatting
$ws1.Range("I4").Value = ""
$ws1.Range("J4").Value = ""
$ws1.Range("K4").Value = ""
$ws1.Range("L4").Value = ""
$ws1.Range("M4").Value = ""

$pctRange = $ws1.Range("I3:M4")
$pctRange.NumberFormat = "0.00%"
$pctRange.HorizontalAlignment = -4108
$pctRange.Borders.LineStyle = 1
$pctRange.Borders.Weight = 2

# Merge the percentage cells vertically (row3:row4) per column
$ws1.Range("I3:I4").Merge()
$ws1.Range("J3:J4").Merge()
$ws1.Range("K3:K4").Merge()
$ws1.Range("L3:L4").Merge()
$ws1.Range("M3:M4").Merge()

# ---------------------------------------------------------------------
# Sheet2 — pass/fail/not run/blocked counts per test-case section
# ---------------------------------------------------------------------

$ws2.Columns.Item(1).ColumnWidth = 11

$ws2.Range("A6").Value = "percentage"
$ws2.Range("B6").Value = "pass"
$ws2.Range("C6").Value = "fail"
$ws2.Range("D6").Value = "not run"
$ws2.Range("E6").Value = "blocked"
$ws2.Range("B6:E6").HorizontalAlignment = -4108

function Add-CountRow($row, $label, $range) {
    $ws2.Range("A$row").Value = $label
    $ws2.Range("B$row").Formula = "=COUNTIF(Sheet1!$range, ""PASS"")"
    $ws2.Range("C$row").Formula = "=COUNTIF(Sheet1!$range, ""FAIL"")"
    $ws2.Range("D$row").Formula = "=COUNTIF(Sheet1!$range, ""NOT RUN"")"
    $ws2.Range("E$row").Formula = "=COUNTIF(Sheet1!$range, ""BLOCKED"")"
    $ws2.Range("A$row`:E$row").HorizontalAlignment = -4108
}

Add-CountRow 7  2     "F7:F13"
Add-CountRow 8  2.1   "F15:F20"
Add-CountRow 9  2.2   "F22:F38"
Add-CountRow 10 2.3   "F40:F52"
Add-CountRow 11 2.4   "F54:F61"
Add-CountRow 12 2.5   "F63:F67"
Add-CountRow 13 "2.5.1" "F69:F74"
Add-CountRow 14 "2.5.2" "F76:F78"
Add-CountRow 15 "2.5.3" "F80:F96"
Add-CountRow 16 "2.5.4" "F98:F100"

# Grand totals
$ws2.Range("A17").Value = "total"
$ws2.Range("B17").Formula = "=SUM(B7:B16)"
$ws2.Range("C17").Formula = "=SUM(C7:C16)"
$ws2.Range("D17").Formula = "=SUM(D7:D16)"
$ws2.Range("E17").Formula = "=SUM(E7:E16)"
$ws2.Range("A17").HorizontalAlignment = -4108

Write-Host "summary tables written"
